$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the date/time value on row 75, column A (keeps its existing style) ---
$ws.Range("A75").Value = 45446.2916666667

# --- Append new row 76 with the latest data point from the R script ---

# Column A: date serial — write the value, then copy row 75's formatting
# (number format yyyy-mm-dd hh:mm:ss) onto it so it matches the other date cells
# without inventing a brand-new style entry.
$ws.Range("A76").Value = 45447.6494328704
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)

# Columns B-F: plain numeric values
$ws.Range("B76").Value = 9000
$ws.Range("C76").Value = 6.69999980926514
$ws.Range("D76").Value = 6.46000003814697
$ws.Range("E76").Value = 6.53999996185303
$ws.Range("F76").Value = 6.42000007629395

# Column G: adj_close stored as text in the shared-string table (new unique string)
$ws.Range("G76").NumberFormat = "@"
$ws.Range("G76").Value = "6.42000007629395"
$ws.Range("G76").ClearFormats()

# Column H: ticker text (reuses the existing "PAL.MI" shared string)
$ws.Range("H76").Value = "PAL.MI"
